$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.574.31"
$ws.Range("E2").Value = "  +0.00%  "
$ws.Range("D3").Value = "1.598.63"
$ws.Range("E3").Value = "  -0.33%  "
$ws.Range("E4").Value = "  +1.20%  "
$ws.Range("D5").Value = "'208.58"
$ws.Range("E5").Value = "  +0.61%  "
$ws.Range("D6").Value = "'0.501"
$ws.Range("E6").Value = "  -3.13%  "
$ws.Range("E7").Value = "  +1.61%  "
$ws.Range("D8").Value = "'22.35"
$ws.Range("E8").Value = "  -2.96%  "
$ws.Range("D9").Value = "'0.251"
$ws.Range("E9").Value = "  -0.84%  "
$ws.Range("E10").Value = "  -2.30%  "
$ws.Range("D11").Value = "'0.0871"
$ws.Range("E11").Value = "  +0.29%  "
$ws.Range("D12").Value = "1.824.91"
$ws.Range("E12").Value = "  -0.71%  "
$ws.Range("D13").Value = "1.585.39"
$ws.Range("E13").Value = "  -0.72%  "
$ws.Range("E14").Value = "  -2.81%  "
$ws.Range("D15").Value = "'0.541"
$ws.Range("E15").Value = "  -2.30%  "
$ws.Range("D16").Value = "'63.57"
$ws.Range("E16").Value = "  -1.41%  "
$ws.Range("D17").Value = "27.554.50"
$ws.Range("E17").Value = "  -0.06%  "
$ws.Range("D18").Value = "'216.96"
$ws.Range("E18").Value = "  -4.24%  "
$ws.Range("D19").Value = "'7.42"
$ws.Range("E19").Value = "  -1.07%  "
$ws.Range("D20").Value = "0.0₃0691"
$ws.Range("E20").Value = "  -4.03%  "
$ws.Range("E21").Value = "  +1.26%  "
$ws.Range("E22").Value = "  -1.22%  "
$ws.Range("D23").Value = "'9.78"
$ws.Range("E23").Value = "  -2.06%  "
$ws.Range("E24").Value = "  +0.07%  "
$ws.Range("D25").Value = "'154.43"
$ws.Range("E25").Value = "  +1.10%  "
$ws.Range("E26").Value = "  +1.17%  "
$ws.Range("D27").Value = "'6.71"
$ws.Range("E27").Value = "  -1.67%  "
$ws.Range("D28").Value = "'15.04"
$ws.Range("E28").Value = "  -1.88%  "
$ws.Range("E29").Value = "  -3.74%  "
$ws.Range("E30").Value = "  +0.26%  "
$ws.Range("D31").Value = "'0.0467"
$ws.Range("E31").Value = "  -1.76%  "
$ws.Range("E32").Value = "  -1.85%  "
$ws.Range("D33").Value = "1.367.38"
$ws.Range("E33").Value = "  -0.22%  "
$ws.Range("E34").Value = "  -2.43%  "
$ws.Range("E35").Value = "  +0.31%  "
$ws.Range("D36").Value = "'0.965"
$ws.Range("E36").Value = "  -2.04%  "
$ws.Range("E37").Value = "  +0.18%  "
$ws.Range("E38").Value = "  -1.21%  "
$ws.Range("D39").Value = "'0.538"
$ws.Range("E39").Value = "  -1.87%  "
$ws.Range("D40").Value = "'0.816"
$ws.Range("E40").Value = "  -3.01%  "
$ws.Range("E41").Value = "  +1.43%  "
$ws.Range("D42").Value = "'0.966"
$ws.Range("E42").Value = "  -4.92%  "
$ws.Range("E43").Value = "  -1.11%  "
$ws.Range("D44").Value = "'64.10"
$ws.Range("E44").Value = "  -1.25%  "
$ws.Range("E45").Value = "  -1.91%  "
$ws.Range("D46").Value = "1.735.32"
$ws.Range("E46").Value = "  -1.11%  "
$ws.Range("E47").Value = "  -4.79%  "
$ws.Range("E48").Value = "  +0.81%  "
$ws.Range("D49").Value = "0.0₆01000"
$ws.Range("E49").Value = "  +4.45%  "
$ws.Range("D50").Value = "'0.0970"
$ws.Range("E50").Value = "  -3.05%  "
$ws.Range("E51").Value = "  -0.60%  "
